$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. Values must remain plain text (matching the
# original inline-string cells), so each cell is forced to Text format before
# assignment and then returned to the default (Normal) style afterwards so no
# stray formatting is introduced.
$updates = @{
    "E2" = "-4.00%"
    "D3" = "30.84"
    "E3" = "-4.03%"
    "D4" = "4.886"
    "E4" = "-2.03%"
    "D5" = "0.07143"
    "E5" = "-9.54%"
    "D6" = "1.871"
    "E6" = "-10.79%"
    "D7" = "7.637"
    "E7" = "-2.04%"
    "D8" = "3.760"
    "E8" = "-1.63%"
    "D9" = "0.8959"
    "E9" = "-3.57%"
    "D10" = "0.1645"
    "E10" = "-6.07%"
    "D11" = "0.07528"
    "E11" = "-5.70%"
    "D12" = "0.08150"
    "E12" = "-4.90%"
    "D13" = "0.03001"
    "E13" = "-4.86%"
    "D14" = "0.09994"
    "E14" = "-0.14%"
    "D15" = "0.001491"
    "E15" = "-1.23%"
    "D16" = "0.005841"
    "E16" = "-1.74%"
    "E18" = "-0.13%"
    "D19" = "2.106"
    "E19" = "-7.47%"
    "E21" = "-1.30%"
    "D22" = "4.274"
    "E22" = "-0.08%"
    "D23" = "0.2003"
    "E23" = "11.83%"
    "E24" = "-2.77%"
    "D25" = "0.001214"
    "E25" = "-2.01%"
    "D26" = "0.004657"
    "E26" = "4.46%"
    "E27" = "0.09%"
    "D39" = "0.01638"
    "E39" = "-4.30%"
    "D40" = "0.04345"
    "E40" = "-8.93%"
    "D41" = "0.007378"
    "E41" = "-0.95%"
    "D42" = "0.1307"
    "E42" = "-3.86%"
    "D43" = "0.002006"
    "E43" = "-13.58%"
    "D44" = "0.01015"
    "E44" = "-0.89%"
    "D45" = "0.00005863"
    "E45" = "-2.13%"
    "E46" = "0.09%"
    "D47" = "2.205"
    "E47" = "167.75%"
    "E48" = "-11.47%"
    "D49" = "0.00002103"
    "E49" = "0.09%"
    "E50" = "0.09%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
